# Applies the "bsearch rotated search and min" commit:
#  - Adds two new Binary-Search problems (Koko Eating Bananas,
#    Find Minimum in Rotated Sorted Array) to the "Problems" sheet.
#  - Updates row heights / column widths / selection to match.
#  - Bumps the row height of the Binary Search row on "Algorithms".

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------
# Sheet "Problems": append rows 4 and 5
# ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Problems")

$kokoApproach = @'
1. Initialize l, r to 1, max(piles) and res to r or max(piles)
2. While l <= r:
	initialize hours = 0 
	k to mid i.e (l+r)//2
	iterate on p of piles:
		add math.ceil(p / k) to hours
	if hours <= h:
		init res = min (res, k)
		shift r to low range		
	else:
		shift l to high range
3. return res
'@

$findMinApproach = @'
1. initialize result to any random element
2. initialize l and r
3. While l<=r
	if nums[l] < nums[r] that means arr is sorted assign res = minimum of res, nums[l] and break
	inititalize mid and assign res = min(res, nums[m]) (# at some point m and l will be equal)
	if l <= m, move l pointer to the right
	else move, move r pointer to the left
4. return result
'@

# Row 4 - Koko Eating Bananas
$ws.Cells.Item(4, 1).Value2 = "Array"
$ws.Cells.Item(4, 2).Value2 = "Binary Search"
$ws.Cells.Item(4, 3).Value2 = "Koko Eating Bananas"
$ws.Cells.Item(4, 4).Value2 = $kokoApproach
$ws.Cells.Item(4, 4).WrapText = $true
$ws.Rows.Item(4).RowHeight = 195

# Row 5 - Find Minimum in Rotated Sorted Array
$ws.Cells.Item(5, 1).Value2 = "Array"
$ws.Cells.Item(5, 2).Value2 = "Binary Search"
$ws.Cells.Item(5, 3).Value2 = "Find Minimum in Rotated Sorted Array"
$ws.Cells.Item(5, 4).Value2 = $findMinApproach
$ws.Cells.Item(5, 4).WrapText = $true
$ws.Rows.Item(5).RowHeight = 150

# Column widths (best-fit recalculated by Excel once the longer strings
# were added to columns B, C & D)
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Columns.Item(2).ColumnWidth = 34
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 74.6666666666667

# Selection / scroll position left by the author after typing the new rows
$ws.Range("D8").Select()

# ----------------------------------------------------------------
# Sheet "Algorithms": the Binary Search implementation row grew taller
# ----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Algorithms")
$ws2.Rows.Item(2).RowHeight = 105
$ws2.Columns.Item(1).ColumnWidth = 11.6666666666667
$ws2.Columns.Item(2).ColumnWidth = 55.6666666666667
